# Add data for 2022-03-23 (updates the "through" date from March 14 to March 15
# and refreshes the year-over-year carjacking counts for the matching date).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the worksheet and update the column header / shared string text to
# reflect the new "through" date.
$ws.Name = "Through 2022-03-15"
$ws.Range("B1").Value = "March 2022 (through March 15)"

# Update existing counts that changed.
$ws.Range("K3").Value = 4    # Austin / March 2019
$ws.Range("E4").Value = 6    # North Lawndale / March 2021
$ws.Range("T4").Value = 2    # North Lawndale / March 2016
$ws.Range("Q5").Value = 7    # Garfield Park / March 2017
$ws.Range("H7").Value = 2    # South Shore / March 2020
$ws.Range("Q7").Value = 2    # South Shore / March 2017
$ws.Range("N11").Value = 4   # Englewood / March 2018
$ws.Range("T11").Value = 3   # Englewood / March 2016
$ws.Range("B15").Value = 3   # Humboldt Park / March 2022 (through March 15)

# Fill in counts for cells that previously had no data.
$ws.Range("W7").Value = 1    # South Shore / March 2015
$ws.Range("Q26").Value = 1   # Grand Crossing / March 2017
$ws.Range("E29").Value = 1   # Lower West Side / March 2021
$ws.Range("N34").Value = 1   # River North / March 2018
$ws.Range("T36").Value = 1   # Roseland / March 2016
$ws.Range("K48").Value = 1   # Avondale / March 2019
$ws.Range("N49").Value = 1   # Little Village / March 2018

Write-Host "Applied carjacking data update for 2022-03-23"
